$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "86.856.15"
$ws.Range("E2").Value = "  +6.97%  "

$ws.Range("D3").Value = "3.305.88"
$ws.Range("E3").Value = "  +2.76%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'217.69"
$ws.Range("E5").Value = "  +2.15%  "

$ws.Range("D6").Value = "'632.58"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").Value = "'0.325"
$ws.Range("E7").Value = "  +14.48%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.610"
$ws.Range("E9").Value = "  +2.97%  "

$ws.Range("D10").Value = "3.299.62"
$ws.Range("E10").Value = "  +2.68%  "

$ws.Range("D11").Value = "'0.597"
$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("D12").Value = "'0.0000270"
$ws.Range("E12").Value = "  +2.36%  "

$ws.Range("E13").Value = "  +1.55%  "

$ws.Range("D14").Value = "3.913.02"
$ws.Range("E14").Value = "  +3.18%  "

$ws.Range("D15").Value = "'34.20"
$ws.Range("E15").Value = "  +5.84%  "

$ws.Range("D16").Value = "'5.36"
$ws.Range("E16").Value = "  +1.06%  "

$ws.Range("D17").Value = "86.776.92"
$ws.Range("E17").Value = "  +7.45%  "

$ws.Range("D18").Value = "3.349.01"
$ws.Range("E18").Value = "  +5.60%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'14.39"
$ws.Range("E19").Value = "  -0.35%  "

$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").Value = "'3.16"
$ws.Range("E20").Value = "  +3.71%  "

$ws.Range("D21").Value = "'444.70"
$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("D22").Value = "'8.93"
$ws.Range("E22").Value = "  -3.91%  "

$ws.Range("D23").Value = "'5.27"
$ws.Range("E23").Value = "  +0.81%  "

$ws.Range("D24").Value = "'7.36"
$ws.Range("E24").Value = "  +5.61%  "

$ws.Range("D25").Value = "'5.32"
$ws.Range("E25").Value = "  +12.29%  "

$ws.Range("D26").Value = "'12.29"
$ws.Range("E26").Value = "  +11.59%  "

$ws.Range("D27").Value = "3.507.82"
$ws.Range("E27").Value = "  +5.18%  "

$ws.Range("D28").Value = "'77.77"
$ws.Range("E28").Value = "  +1.36%  "

$ws.Range("E29").Value = "  +72.31%  "

$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -0.13%  "

$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "'0.0000127"
$ws.Range("E31").Value = "  +1.96%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'9.17"
$ws.Range("E32").Value = "  +1.21%  "

$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'588.90"
$ws.Range("E33").Value = "  +3.02%  "

$ws.Range("E34").Value = "  +0.68%  "

$ws.Range("D35").Value = "'1.51"
$ws.Range("E35").Value = "  +1.45%  "

$ws.Range("D36").Value = "'2.03"
$ws.Range("E36").Value = "  +1.09%  "

$ws.Range("D37").Value = "'0.151"
$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("D38").Value = "'23.15"
$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("D39").Value = "'6.56"
$ws.Range("E39").Value = "  +14.18%  "

$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("D41").Value = "'0.413"
$ws.Range("E41").Value = "  +1.08%  "

$ws.Range("E42").Value = "  +2.84%  "

$ws.Range("D43").Value = "'2.04"
$ws.Range("E43").Value = "  +11.43%  "

$ws.Range("D44").Value = "'3.03"
$ws.Range("E44").Value = "  +10.37%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'158.23"
$ws.Range("E45").Value = "  -3.89%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").Value = "'187.17"
$ws.Range("E47").Value = "  -2.79%  "

$ws.Range("D48").Value = "'46.27"
$ws.Range("E48").Value = "  +7.73%  "

$ws.Range("D49").Value = "'1.35"
$ws.Range("E49").Value = "  +2.87%  "

$ws.Range("D50").Value = "'0.778"
$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("D51").Value = "'26.19"
$ws.Range("E51").Value = "  +4.92%  "
